# ==========================================================================
# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets
#
# 1. Inserts a new "Player Info" sheet before the existing "ODI Batting"
#    sheet with the player's biographical data.
# 2. Updates the "ODI Batting" sheet: renames MATCH_CARD_LINK -> MATCH_CODE
#    and replaces the full scorecard URL with the bare match code; also
#    drops the stray empty B44 cell.
# 3. Appends a new "ODI Batting Extra" sheet after "ODI Batting" holding
#    additional per-match batting stats.
# ==========================================================================

$wb = $excel.ActiveWorkbook

# --------------------------------------------------------------------
# 1. "Player Info" sheet - inserted before "ODI Batting"
# --------------------------------------------------------------------
$odi = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($odi)
$playerInfo.Name = "Player Info"

$piHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($col = 1; $col -le $piHeaders.Length; $col++) {
    $cell = $playerInfo.Cells.Item(1, $col)
    $cell.Value = $piHeaders[$col - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$playerInfo.Cells.Item(2, 1).Value = "'4435"
$playerInfo.Cells.Item(2, 2).Value = "Henry George Munsey"
$playerInfo.Cells.Item(2, 3).Value = "Right Handed"
$playerInfo.Cells.Item(2, 4).Value = "Right Arm Medium Fast"

$playerInfo.Range("A1").Select()

# --------------------------------------------------------------------
# 2. "ODI Batting" sheet updates
#    (re-fetch the sheet reference - it can go stale after the
#    worksheet collection was mutated by the Add() above)
# --------------------------------------------------------------------
$odi = $wb.Worksheets.Item("ODI Batting")

$odi.Cells.Item(1, 4).Value = "MATCH_CODE"

for ($r = 2; $r -le 50; $r++) {
    $cell = $odi.Cells.Item($r, 4)
    $link = $cell.Value2
    $parts = $link -split "MatchCode="
    $code = $parts[$parts.Length - 1]
    $cell.Value = "'" + $code
}

# Drop the stray empty cell at B44
$odi.Cells.Item(44, 2).Value = ""

# --------------------------------------------------------------------
# 3. "ODI Batting Extra" sheet - inserted after "ODI Batting"
# --------------------------------------------------------------------
$odi = $wb.Worksheets.Item("ODI Batting")
$extra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $odi)
$extra.Name = "ODI Batting Extra"

$exHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($col = 1; $col -le $exHeaders.Length; $col++) {
    $cell = $extra.Cells.Item(1, $col)
    $cell.Value = $exHeaders[$col - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$exData = @(
    @("4513", $null, $null,  $null, $null,    "NO"),
    @("4515", $null, $null,  $null, $null,    "NO"),
    @("4569", 5,     "6",    "4",   "27.82%", "NO"),
    @("4570", 5,     "6",    "0",   "28.84%", "NO"),
    @("4572", 5,     "3",    "1",   "12.54%", "NO"),
    @("4573", $null, $null,  $null, $null,    "NO"),
    @("4575", 5,     "7",    "1",   "26.21%", "NO"),
    @("4576", 5,     "7",    "0",   "17.00%", "NO"),
    @("4578", 5,     "0",    "0",   "0.46%",  "NO"),
    @("4581", 5,     "1",    "0",   "16.96%", "NO"),
    @("4604", 5,     "7",    "1",   "27.52%", "NO"),
    @("4610", 5,     "0",    "0",   "2.78%",  "NO"),
    @("4612", 5,     "2",    "0",   "7.76%",  "NO"),
    @("4617", 5,     $null,  $null, $null,    "NO"),
    @("4677", 2,     "1",    "0",   "2.35%",  "NO"),
    @("4681", 1,     "3",    "0",   "11.89%", "NO"),
    @("4680", 2,     "6",    "1",   "24.52%", "NO"),
    @("4684", 2,     "7",    "2",   "37.19%", "NO"),
    @("4702", 2,     "12",   "7",   "65.61%", "YES"),
    @("4706", $null, $null,  $null, $null,    "NO")
)

$rowIdx = 2
foreach ($row in $exData) {
    $extra.Cells.Item($rowIdx, 1).Value = "'" + $row[0]

    if ($null -ne $row[1]) {
        $extra.Cells.Item($rowIdx, 2).Value = $row[1]
    }
    if ($null -ne $row[2]) {
        $extra.Cells.Item($rowIdx, 3).Value = "'" + $row[2]
    }
    if ($null -ne $row[3]) {
        $extra.Cells.Item($rowIdx, 4).Value = "'" + $row[3]
    }
    if ($null -ne $row[4]) {
        $extra.Cells.Item($rowIdx, 5).Value = "'" + $row[4]
    }
    $extra.Cells.Item($rowIdx, 6).Value = $row[5]

    $rowIdx++
}

$extra.Range("A1").Select()
$odi.Activate()
